# Generate Report for Handback
# Updates the localization-status workbook to reflect a handback event:
#  - "Ready for handoff" -> "Handed back: in sync with en-US" status text
#  - zh-cn handback datetime refreshed (2016-09-03 02:32:37)
#  - de-de handback datetime recorded for the first time (2016-09-03 02:32:44)
#  - Latest Target File / Latest Handback File columns populated + hyperlinked
#    for both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet: status moves from "Ready for handoff" to
# "Handed back: in sync with en-US" for both tracked files/languages.
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: record the target (.md) file and handback (.xlf) file for
# each source row, and refresh the handback datetime.
# ---------------------------------------------------------------------------
$wsZh.Range("I2").Value = "22ed719f-21c2-4705-85b3-513b44286fe1.md"
$wsZh.Range("J2").Value = "22ed719f-21c2-4705-85b3-513b44286fe1.e8e5b062479d10817b2a19624f5ebe94b0a1cc9a.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 02:32:37"

$wsZh.Range("I3").Value = "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md"
$wsZh.Range("J3").Value = "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.2c0fb9a989ace7f33443fdd4d061578bed07bb62.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 02:32:37"

# Rebuild the hyperlinks so column I picks up the same link as column A,
# keeping rId ordering in row-major order (A2, I2, A3, I3).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/22ed719f-21c2-4705-85b3-513b44286fe1.md", [System.Type]::Missing, [System.Type]::Missing, "22ed719f-21c2-4705-85b3-513b44286fe1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/22ed719f-21c2-4705-85b3-513b44286fe1.md", [System.Type]::Missing, [System.Type]::Missing, "22ed719f-21c2-4705-85b3-513b44286fe1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md", [System.Type]::Missing, [System.Type]::Missing, "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md", [System.Type]::Missing, [System.Type]::Missing, "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md")

# ---------------------------------------------------------------------------
# de-de sheet: same shape of update, but this is the *first* handback, so
# the datetime is brand new (not a refresh of an existing value).
# ---------------------------------------------------------------------------
$wsDe.Range("I2").Value = "22ed719f-21c2-4705-85b3-513b44286fe1.md"
$wsDe.Range("J2").Value = "22ed719f-21c2-4705-85b3-513b44286fe1.e8e5b062479d10817b2a19624f5ebe94b0a1cc9a.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 02:32:44"

$wsDe.Range("I3").Value = "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md"
$wsDe.Range("J3").Value = "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.2c0fb9a989ace7f33443fdd4d061578bed07bb62.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 02:32:44"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/22ed719f-21c2-4705-85b3-513b44286fe1.md", [System.Type]::Missing, [System.Type]::Missing, "22ed719f-21c2-4705-85b3-513b44286fe1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/22ed719f-21c2-4705-85b3-513b44286fe1.md", [System.Type]::Missing, [System.Type]::Missing, "22ed719f-21c2-4705-85b3-513b44286fe1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md", [System.Type]::Missing, [System.Type]::Missing, "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/389c0d4df0d8d7a361da982fbe18d1f9959cb71a/e2e/6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md", [System.Type]::Missing, [System.Type]::Missing, "6f54e31e-ed79-4bf6-8b34-e22f4c8fad48.md")

# ---------------------------------------------------------------------------
# Column widths: the newly populated / widened columns auto-fit wider once
# the handback file names are present.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(9).ColumnWidth = 39.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(9).ColumnWidth = 39.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15

Write-Output "handback report generated"
